$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 7
$ws.Range("I8").Value = 7
$ws.Range("K8").Value = 21
$ws.Range("M8").Value = 118

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 524.25
$ws.Range("I9").Value = 568.375
$ws.Range("J9").Value = 436
$ws.Range("K9").Value = 568.375
$ws.Range("L9").Value = 436
$ws.Range("M9").Value = -399.375
$ws.Range("N9").Value = -774

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 10870642
$ws.Range("I15").Value = 10870642
$ws.Range("K15").Value = 32611926
$ws.Range("M15").Value = -32611757

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1201.2766
$ws.Range("J17").Value = 1201.2766
$ws.Range("L17").Value = 3603.8298
$ws.Range("N17").Value = -3939.8298

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 876
$ws.Range("J100").Value = 699
$ws.Range("L100").Value = 699
$ws.Range("N100").Value = -1781

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 3792.12
$ws.Range("J121").Value = 3792.12
$ws.Range("L121").Value = 11376.36
$ws.Range("N121").Value = -14870.36

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2597
$ws.Range("I127").Value = 2422.3333
$ws.Range("K127").Value = 7266.999899999999
$ws.Range("M127").Value = -2306.999899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1877.5483
$ws.Range("I135").Value = 381.45834
$ws.Range("K135").Value = 3433.12506
$ws.Range("M135").Value = -898.1250600000003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9808732
$ws.Range("I137").Value = 1868.3334
$ws.Range("K137").Value = 5605.0002
$ws.Range("M137").Value = -3055.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3257.1667
$ws.Range("I138").Value = 2088.5
$ws.Range("J138").Value = 3424.1191
$ws.Range("K138").Value = 6265.5
$ws.Range("L138").Value = 10272.3573
$ws.Range("M138").Value = -1125.5
$ws.Range("N138").Value = -20552.3573

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 15582.143
$ws.Range("I141").Value = 15582.143
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 46746.429
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -41566.429
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4328.537
$ws.Range("I32").Value = 2595.889
$ws.Range("J32").Value = 11525.692
$ws.Range("K32").Value = 2595.889
$ws.Range("L32").Value = 11525.692
$ws.Range("M32").Value = -2308.889
$ws.Range("N32").Value = -12099.692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 5281
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 5559.1
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 5559.1
$ws.Range("M46").Value = -2181
$ws.Range("N46").Value = -6197.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14151.25
$ws.Range("I61").Value = 15256
$ws.Range("K61").Value = 15256
$ws.Range("M61").Value = -15044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13159840
$ws.Range("I74").Value = 19232252
$ws.Range("K74").Value = 19232252
$ws.Range("M74").Value = -19231378

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 13159840
$ws.Range("I77").Value = 19232252
$ws.Range("K77").Value = 96161260
$ws.Range("M77").Value = -96156892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4536.533
$ws.Range("I122").Value = 3049.6667
$ws.Range("K122").Value = 9149.000100000001
$ws.Range("M122").Value = -6699.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 14151.25
$ws.Range("I136").Value = 15256
$ws.Range("K136").Value = 45768
$ws.Range("M136").Value = -43218

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2084835.8
$ws.Range("I99").Value = 2977479.8
$ws.Range("J99").Value = 1999.6666
$ws.Range("K99").Value = 2977479.8
$ws.Range("L99").Value = 1999.6666
$ws.Range("M99").Value = -2975981.8
$ws.Range("N99").Value = -4995.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24396702
$ws.Range("I31").Value = 100001600
$ws.Range("J31").Value = 8026.032
$ws.Range("K31").Value = 100001600
$ws.Range("L31").Value = 8026.032
$ws.Range("M31").Value = -100001305
$ws.Range("N31").Value = -8616.031999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 24396702
$ws.Range("I34").Value = 100001600
$ws.Range("J34").Value = 8026.032
$ws.Range("K34").Value = 100001600
$ws.Range("L34").Value = 8026.032
$ws.Range("M34").Value = -100001398
$ws.Range("N34").Value = -8430.031999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 102597010
$ws.Range("I132").Value = 133355620
$ws.Range("K132").Value = 400066860
$ws.Range("M132").Value = -400064330

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4090
$ws.Range("I134").Value = 4090
$ws.Range("K134").Value = 12270
$ws.Range("M134").Value = -9735

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 257.16666
$ws.Range("I14").Value = 257.16666
$ws.Range("K14").Value = 771.4999799999999
$ws.Range("M14").Value = -598.4999799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 480
$ws.Range("I44").Value = 306.66666
$ws.Range("J44").Value = 1000
$ws.Range("K44").Value = 919.9999799999999
$ws.Range("L44").Value = 3000
$ws.Range("M44").Value = -521.9999799999999
$ws.Range("N44").Value = -3796

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2307.2856
$ws.Range("J107").Value = 2307.2856
$ws.Range("L107").Value = 6921.8568
$ws.Range("N107").Value = -10761.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 5598.75
$ws.Range("I118").Value = 2698.5
$ws.Range("J118").Value = 8499
$ws.Range("K118").Value = 8095.5
$ws.Range("L118").Value = 25497
$ws.Range("M118").Value = -6852.5
$ws.Range("N118").Value = -27983

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6452331.5
$ws.Range("I122").Value = 32258064
$ws.Range("J122").Value = 898.5
$ws.Range("K122").Value = 290322576
$ws.Range("L122").Value = 8086.5
$ws.Range("M122").Value = -290320126
$ws.Range("N122").Value = -12986.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3695.889
$ws.Range("J140").Value = 5189.5713
$ws.Range("L140").Value = 15568.7139
$ws.Range("N140").Value = -25928.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1592337.4
$ws.Range("J70").Value = 5687.125
$ws.Range("L70").Value = 5687.125
$ws.Range("N70").Value = -6227.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 1592337.4
$ws.Range("J73").Value = 5687.125
$ws.Range("L73").Value = 5687.125
$ws.Range("N73").Value = -7559.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 25008896
$ws.Range("I102").Value = 35723696
$ws.Range("K102").Value = 35723696
$ws.Range("M102").Value = -35722074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1176.375
$ws.Range("I113").Value = 485.16666
$ws.Range("J113").Value = 3250
$ws.Range("K113").Value = 485.16666
$ws.Range("L113").Value = 3250
$ws.Range("M113").Value = 1684.83334
$ws.Range("N113").Value = -7590

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 462020.38
$ws.Range("I122").Value = 613833
$ws.Range("K122").Value = 1841499
$ws.Range("M122").Value = -1839049

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 44291.617
$ws.Range("J123").Value = 44291.617
$ws.Range("L123").Value = 44291.617
$ws.Range("N123").Value = -49191.617

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3143.6572
$ws.Range("I132").Value = 2507.889
$ws.Range("K132").Value = 7523.667
$ws.Range("M132").Value = -4993.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 34833
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 49999.5
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 49999.5
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -51247.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 34833
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 49999.5
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 249997.5
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -256237.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 150000
$ws.Range("J80").Value = 150000
$ws.Range("L80").Value = 150000
$ws.Range("N80").Value = -151996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 150000
$ws.Range("J83").Value = 150000
$ws.Range("L83").Value = 450000
$ws.Range("N83").Value = -459984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 728.6957
$ws.Range("I113").Value = 690
$ws.Range("J113").Value = 788.8889
$ws.Range("K113").Value = 2070
$ws.Range("L113").Value = 2366.6667
$ws.Range("M113").Value = 100
$ws.Range("N113").Value = -6706.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5129.5947
$ws.Range("I122").Value = 4446.1113
$ws.Range("J122").Value = 6975
$ws.Range("K122").Value = 13338.3339
$ws.Range("L122").Value = 20925
$ws.Range("M122").Value = -10888.3339
$ws.Range("N122").Value = -25825

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22227262
$ws.Range("I132").Value = 3088087.5
$ws.Range("K132").Value = 9264262.5
$ws.Range("M132").Value = -9261732.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9236.986000000001
$ws.Range("I136").Value = 5216.643
$ws.Range("J136").Value = 11795.387
$ws.Range("K136").Value = 15649.929
$ws.Range("L136").Value = 35386.161
$ws.Range("M136").Value = -13099.929
$ws.Range("N136").Value = -40386.161
